# Generate Report for Handback
# Refresh the handback/handoff timestamps recorded on the "Overview", "zh-cn"
# and "de-de" sheets (row 2 = c4ced760-1093-4bf5-b9c5-270ed45bf59a.md).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-26 03:04:31"

# zh-cn!H2 - Correspond Handoff Datetime
$wsZhCn.Range("H2").Value = "2016-08-26 03:04:26"

# zh-cn!K2 - Correspond Handback DateTime
$wsZhCn.Range("K2").Value = "2016-08-26 03:04:48"

# de-de!H2 - Correspond Handoff Datetime (shares text with Overview!G2)
$wsDeDe.Range("H2").Value = "2016-08-26 03:04:31"

# de-de!K2 - Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-26 03:04:54"
